$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.564.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.988.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.46%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '497.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("E10").Value = '  -2.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.354'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.88%  '
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.494.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.504.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.985.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.73%  '
$ws.Range("E17").Value = '  -3.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.86'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.468'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.12%  '
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.163'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0{0}0902' -f [char]0x2083
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.45%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.77%  '
$ws.Range("E31").Value = '  -3.04%  '
$ws.Range("E32").Value = '  -6.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '155.68'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0685'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.016.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.638'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.31%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.995'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.25%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.226.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.91%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.00%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0236'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.17%  '
